$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("美住登录")

# The "实际结果" (actual result) column D was filled with a redundant
# "PASS" literal on every data row (D2:D9). Clear that stale expected
# content out, leaving the D1 header ("实际结果") untouched.
$ws.Range("D2:D9").ClearContents() | Out-Null

# Leave the cursor on the last touched cell, matching the saved selection.
$ws.Range("D9").Select() | Out-Null
